$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.898999999999999
$ws.Range("B3").Value = 6.265
$ws.Range("D3").Value = -7.388
$ws.Range("B4").Value = 7.092999999999999
$ws.Range("E8").Value = 16.789
$ws.Range("D9").Value = -7.232000000000001
$ws.Range("A11").Value = -21.413
$ws.Range("E11").Value = 17.139
$ws.Range("A12").Value = -21.609
$ws.Range("B14").Value = 6.199999999999999
$ws.Range("E14").Value = 17.003
$ws.Range("A15").Value = -21.202
$ws.Range("D15").Value = -7.938000000000001
$ws.Range("E15").Value = 16.576
$ws.Range("E17").Value = 16.725
$ws.Range("D19").Value = -8.159000000000001
$ws.Range("D20").Value = -7.765000000000001
$ws.Range("D25").Value = -7.938
$ws.Range("B26").Value = 6.459000000000001
$ws.Range("E26").Value = 17.027
$ws.Range("A27").Value = -21.026
$ws.Range("D27").Value = -8.247
$ws.Range("A28").Value = -21.527
$ws.Range("D28").Value = -8.028
$ws.Range("D30").Value = -7.322
$ws.Range("A31").Value = -21.29
$ws.Range("B31").Value = 6.223000000000001
$ws.Range("A32").Value = -21.234
$ws.Range("D32").Value = -8.087
$ws.Range("B35").Value = 7.216000000000001
$ws.Range("A36").Value = -21.156
$ws.Range("E36").Value = 16.767
$ws.Range("B37").Value = 7.065
$ws.Range("A38").Value = -19.741
$ws.Range("B39").Value = 7.615
$ws.Range("B40").Value = 8.211
$ws.Range("E42").Value = 16.494
$ws.Range("D44").Value = -7.956
$ws.Range("B45").Value = 6.226000000000001
$ws.Range("A46").Value = -21.266
$ws.Range("D47").Value = -7.489
$ws.Range("B52").Value = 4.883000000000001
$ws.Range("A54").Value = -21.856
$ws.Range("A55").Value = -22.21
$ws.Range("A56").Value = -21.803
$ws.Range("B57").Value = 5.331999999999999
$ws.Range("D58").Value = -8.146000000000001
$ws.Range("D62").Value = -8.087
$ws.Range("E64").Value = 17.145
$ws.Range("A67").Value = -21.603
$ws.Range("E68").Value = 17.301
$ws.Range("A69").Value = -21.636
$ws.Range("A72").Value = -21.567
$ws.Range("A73").Value = -20.71
$ws.Range("D77").Value = -7.891000000000001
$ws.Range("D78").Value = -7.812
$ws.Range("E79").Value = 17.266
$ws.Range("B81").Value = 6.392000000000001
$ws.Range("A83").Value = -20.531
$ws.Range("B83").Value = 6.208
$ws.Range("D84").Value = -8.148
$ws.Range("A86").Value = -21.962
$ws.Range("D89").Value = -7.105
$ws.Range("E89").Value = 17.352
$ws.Range("A91").Value = -21.534
$ws.Range("D91").Value = -7.139999999999999
$ws.Range("D92").Value = -7.056999999999999
$ws.Range("A93").Value = -21.665
$ws.Range("D96").Value = -7.606999999999999
$ws.Range("A99").Value = -20.547
$ws.Range("B100").Value = 5.558
$ws.Range("B102").Value = 7.499000000000001
$ws.Range("D102").Value = -7.986999999999999
